# Update gh-pages to output generated at 456a3b4
# Applies the numeric "want-to-go" count bumps (column F) and marks one
# event as sold out ("不可售" in column G) on both the "展览" sheet and
# the mirrored "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F4").Value = 273
        $ws.Range("F6").Value = 57
        $ws.Range("F8").Value = 215
        $ws.Range("G8").Value = "不可售"
        $ws.Range("F9").Value = 1988
        $ws.Range("F11").Value = 4677
        $ws.Range("F12").Value = 86
    }
    else {
        $ws.Range("F6").Value = 273
        $ws.Range("F8").Value = 57
        $ws.Range("F10").Value = 215
        $ws.Range("G10").Value = "不可售"
        $ws.Range("F13").Value = 1988
        $ws.Range("F15").Value = 4677
        $ws.Range("F16").Value = 86
    }
}

Write-Host "Applied updates to sheets: $($sheetNames -join ', ')"
